$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(47, 1).Value = "SXT"
$ws.Cells.Item(47, 2).Value = "Year"
$ws.Cells.Item(47, 3).Value = 0.001920012610919285

$ws.Cells.Item(48, 1).Value = "SXT"
$ws.Cells.Item(48, 2).Value = "Specimen_type"
$ws.Cells.Item(48, 3).Value = 0.06206679387812179

$ws.Cells.Item(49, 1).Value = "SXT"
$ws.Cells.Item(49, 2).Value = "Gender"
$ws.Cells.Item(49, 3).Value = 0.00008739521480219351

$ws.Cells.Item(50, 1).Value = "SXT"
$ws.Cells.Item(50, 2).Value = "Age_cat"
$ws.Cells.Item(50, 3).Value = 0.4722429984634713

$ws.Cells.Item(51, 1).Value = "SXT"
$ws.Cells.Item(51, 2).Value = "Hospital:Ward_ED_ICU"
$ws.Cells.Item(51, 3).Value = 0.0005552045233266024
